$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Rows with Price (D) and Volume(1h) (E) updates
Set-TextValue "D2" "67.162.81"
$ws.Range("E2").Value = "  +1.52%  "
Set-TextValue "D3" "3.858.19"
$ws.Range("E3").Value = "  +0.99%  "
Set-TextValue "D5" "467.80"
$ws.Range("E5").Value = "  +9.76%  "
Set-TextValue "D6" "144.68"
$ws.Range("E6").Value = "  +10.70%  "
Set-TextValue "D7" "0.632"
$ws.Range("E7").Value = "  +3.47%  "
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.745"
$ws.Range("E9").Value = "  +2.30%  "
Set-TextValue "D11" "0.0000310"
$ws.Range("E11").Value = "  -6.82%  "
Set-TextValue "D12" "43.42"
$ws.Range("E12").Value = "  +4.55%  "
Set-TextValue "D13" "10.44"
$ws.Range("E13").Value = "  -0.08%  "
Set-TextValue "D14" "4.488.77"
$ws.Range("E14").Value = "  +1.12%  "
Set-TextValue "D15" "14.83"
$ws.Range("E15").Value = "  -4.55%  "
Set-TextValue "D16" "3.904.26"
$ws.Range("E16").Value = "  +0.40%  "
Set-TextValue "D19" "1.16"
$ws.Range("E19").Value = "  +5.86%  "
Set-TextValue "D20" "67.424.47"
$ws.Range("E20").Value = "  +1.44%  "
Set-TextValue "D21" "434.94"
$ws.Range("E21").Value = "  +4.69%  "
Set-TextValue "D22" "14.92"
$ws.Range("E22").Value = "  -1.00%  "
Set-TextValue "D23" "3.32"
$ws.Range("E23").Value = "  +6.64%  "
Set-TextValue "D24" "88.92"
$ws.Range("E24").Value = "  +4.67%  "
Set-TextValue "D25" "3.60"
$ws.Range("E25").Value = "  +9.51%  "
Set-TextValue "D26" "37.92"
$ws.Range("E26").Value = "  +1.87%  "
Set-TextValue "D27" "10.13"
$ws.Range("E27").Value = "  +7.28%  "
Set-TextValue "D29" "5.54"
$ws.Range("E29").Value = "  +2.77%  "
Set-TextValue "D30" "728.14"
$ws.Range("E30").Value = "  +1.44%  "
Set-TextValue "D33" "2.79"
$ws.Range("E33").Value = "  +3.85%  "
Set-TextValue "D34" "44.21"
$ws.Range("E34").Value = "  +13.47%  "
Set-TextValue "D35" "0.160"
$ws.Range("E35").Value = "  +7.05%  "
Set-TextValue "D36" "58.16"
$ws.Range("E36").Value = "  +4.50%  "
Set-TextValue "D38" "5.49"
$ws.Range("E38").Value = "  -3.80%  "
Set-TextValue "D43" "0.0₃0675"
$ws.Range("E43").Value = "  -7.24%  "
Set-TextValue "D45" "2.57"
$ws.Range("E45").Value = "  +7.84%  "
Set-TextValue "D46" "3.45"
$ws.Range("E46").Value = "  +1.95%  "
Set-TextValue "D47" "3.27"
$ws.Range("E47").Value = "  -0.18%  "
Set-TextValue "D48" "2.78"
$ws.Range("E48").Value = "  +5.58%  "

# Rows with only Volume(1h) (E) updates
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  +6.90%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +4.48%  "

# Rows with full Coin/Link/Price/Volume swap updates (rows reordered in source)
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D40" "0.348"
$ws.Range("E40").Value = "  +8.63%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D41" "2.92"
$ws.Range("E41").Value = "  +1.39%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D50" "2.90"
$ws.Range("E50").Value = "  +1.83%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "144.14"
$ws.Range("E51").Value = "  +1.69%  "

